# Re-price a batch of Leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# sheets with refreshed Universalis market-board averages (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ) and the LeveProfitNQ/HQ figures that
# derive from them (columns H-N). Mirrors the scheduled-runner commit that
# refreshed these quotes; a handful of rows also gain/lose an HQ-profit cell
# depending on whether an HQ market sample now exists.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 177.8
$ws.Cells.Item(9, 9).Value = 139.71428
$ws.Cells.Item(9, 11).Value = 139.71428
$ws.Cells.Item(9, 13).Value = 29.28572

$ws.Cells.Item(11, 8).Value = 4034
$ws.Cells.Item(11, 9).Value = 4034
$ws.Cells.Item(11, 11).Value = 4034
$ws.Cells.Item(11, 13).Value = -3894

$ws.Cells.Item(21, 8).Value = 2000
$ws.Cells.Item(21, 9).Value = 2000
$ws.Cells.Item(21, 11).Value = 2000
$ws.Cells.Item(21, 13).Value = -1532

$ws.Cells.Item(23, 8).Value = 2000
$ws.Cells.Item(23, 9).Value = 2000
$ws.Cells.Item(23, 11).Value = 2000
$ws.Cells.Item(23, 13).Value = -1766

$ws.Cells.Item(43, 8).Value = 1015.26666
$ws.Cells.Item(43, 9).Value = 1006.6923
$ws.Cells.Item(43, 11).Value = 1006.6923
$ws.Cells.Item(43, 13).Value = -937.6923

$ws.Cells.Item(53, 8).Value = 4964.2607
$ws.Cells.Item(53, 9).Value = 227.90909
$ws.Cells.Item(53, 10).Value = 9305.916999999999
$ws.Cells.Item(53, 11).Value = 227.90909
$ws.Cells.Item(53, 12).Value = 9305.916999999999
$ws.Cells.Item(53, 13).Value = 409.09091
$ws.Cells.Item(53, 14).Value = -10579.917

$ws.Cells.Item(69, 8).Value = 47630892
$ws.Cells.Item(69, 10).Value = 52639144
$ws.Cells.Item(69, 12).Value = 157917432
$ws.Cells.Item(69, 14).Value = -157919180

$ws.Cells.Item(72, 8).Value = 47630892
$ws.Cells.Item(72, 10).Value = 52639144
$ws.Cells.Item(72, 12).Value = 473752296
$ws.Cells.Item(72, 14).Value = -473761032

$ws.Cells.Item(116, 8).Value = 5221.357
$ws.Cells.Item(116, 9).Value = 4553.5713
$ws.Cells.Item(116, 10).Value = 5889.143
$ws.Cells.Item(116, 11).Value = 4553.5713
$ws.Cells.Item(116, 12).Value = 5889.143
$ws.Cells.Item(116, 13).Value = -1111.5713
$ws.Cells.Item(116, 14).Value = -12773.143

$ws.Cells.Item(131, 8).Value = 4613.2144
$ws.Cells.Item(131, 9).Value = 3825.4546
$ws.Cells.Item(131, 10).Value = 7501.6665
$ws.Cells.Item(131, 11).Value = 11476.3638
$ws.Cells.Item(131, 12).Value = 22504.9995
$ws.Cells.Item(131, 13).Value = -6436.363799999999
$ws.Cells.Item(131, 14).Value = -32584.9995

$ws.Cells.Item(135, 8).Value = 1090.3939
$ws.Cells.Item(135, 9).Value = 1013.3043
$ws.Cells.Item(135, 10).Value = 1267.7
$ws.Cells.Item(135, 11).Value = 9119.7387
$ws.Cells.Item(135, 12).Value = 11409.3
$ws.Cells.Item(135, 13).Value = -6584.7387
$ws.Cells.Item(135, 14).Value = -16479.3

$ws.Cells.Item(138, 8).Value = 2482.8853
$ws.Cells.Item(138, 9).Value = 1120.4166
$ws.Cells.Item(138, 10).Value = 4444.84
$ws.Cells.Item(138, 11).Value = 3361.2498
$ws.Cells.Item(138, 12).Value = 13334.52
$ws.Cells.Item(138, 13).Value = 1778.7502
$ws.Cells.Item(138, 14).Value = -23614.52

$ws.Cells.Item(141, 8).Value = 1953.8695
$ws.Cells.Item(141, 9).Value = 1727.1666
$ws.Cells.Item(141, 10).Value = 2770
$ws.Cells.Item(141, 11).Value = 5181.4998
$ws.Cells.Item(141, 12).Value = 8310
$ws.Cells.Item(141, 13).Value = -1.499799999999595
$ws.Cells.Item(141, 14).Value = -18670

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 1070.2
$ws.Cells.Item(4, 10).Value = 1002
$ws.Cells.Item(4, 12).Value = 1002
$ws.Cells.Item(4, 14).Value = -1234

$ws.Cells.Item(32, 8).Value = 4245.2095
$ws.Cells.Item(32, 9).Value = 3062.457
$ws.Cells.Item(32, 11).Value = 3062.457
$ws.Cells.Item(32, 13).Value = -2775.457

$ws.Cells.Item(41, 8).Value = 3886.25
$ws.Cells.Item(41, 9).Value = 2181.6667
$ws.Cells.Item(41, 10).Value = 9000
$ws.Cells.Item(41, 11).Value = 2181.6667
$ws.Cells.Item(41, 12).Value = 9000
$ws.Cells.Item(41, 13).Value = -1767.6667
$ws.Cells.Item(41, 14).Value = -9828

$ws.Cells.Item(45, 8).Value = 9054123
$ws.Cells.Item(45, 9).Value = 13988647
$ws.Cells.Item(45, 10).Value = 7495.3335
$ws.Cells.Item(45, 11).Value = 13988647
$ws.Cells.Item(45, 12).Value = 7495.3335
$ws.Cells.Item(45, 13).Value = -13988270
$ws.Cells.Item(45, 14).Value = -8249.333500000001

$ws.Cells.Item(60, 8).Value = 11725
$ws.Cells.Item(60, 9).Value = 9250
$ws.Cells.Item(60, 10).Value = 14200
$ws.Cells.Item(60, 11).Value = 9250
$ws.Cells.Item(60, 12).Value = 14200
$ws.Cells.Item(60, 13).Value = -8517
$ws.Cells.Item(60, 14).Value = -15666

$ws.Cells.Item(61, 8).Value = 4124.0835
$ws.Cells.Item(61, 9).Value = 3319
$ws.Cells.Item(61, 10).Value = 4699.143
$ws.Cells.Item(61, 11).Value = 3319
$ws.Cells.Item(61, 12).Value = 4699.143
$ws.Cells.Item(61, 13).Value = -3107
$ws.Cells.Item(61, 14).Value = -5123.143

$ws.Cells.Item(74, 8).Value = 83862.39999999999
$ws.Cells.Item(74, 9).Value = 12672.091
$ws.Cells.Item(74, 11).Value = 12672.091
$ws.Cells.Item(74, 13).Value = -11798.091

$ws.Cells.Item(77, 8).Value = 83862.39999999999
$ws.Cells.Item(77, 9).Value = 12672.091
$ws.Cells.Item(77, 11).Value = 63360.455
$ws.Cells.Item(77, 13).Value = -58992.455

$ws.Cells.Item(132, 8).Value = 2384.4314
$ws.Cells.Item(132, 9).Value = 1568.3928
$ws.Cells.Item(132, 10).Value = 3377.8696
$ws.Cells.Item(132, 11).Value = 4705.178400000001
$ws.Cells.Item(132, 12).Value = 10133.6088
$ws.Cells.Item(132, 13).Value = -2175.178400000001
$ws.Cells.Item(132, 14).Value = -15193.6088

$ws.Cells.Item(136, 8).Value = 4124.0835
$ws.Cells.Item(136, 9).Value = 3319
$ws.Cells.Item(136, 10).Value = 4699.143
$ws.Cells.Item(136, 11).Value = 9957
$ws.Cells.Item(136, 12).Value = 14097.429
$ws.Cells.Item(136, 13).Value = -7407
$ws.Cells.Item(136, 14).Value = -19197.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 6257690
$ws.Cells.Item(86, 9).Value = 8342947.5
$ws.Cells.Item(86, 11).Value = 8342947.5
$ws.Cells.Item(86, 13).Value = -8341824.5

$ws.Cells.Item(89, 8).Value = 6257690
$ws.Cells.Item(89, 9).Value = 8342947.5
$ws.Cells.Item(89, 11).Value = 41714737.5
$ws.Cells.Item(89, 13).Value = -41709121.5

$ws.Cells.Item(105, 8).Value = 7814598.5
$ws.Cells.Item(105, 9).Value = 7814598.5
$ws.Cells.Item(105, 11).Value = 7814598.5
$ws.Cells.Item(105, 13).Value = -7812851.5

$ws.Cells.Item(107, 8).Value = 7939040.5
$ws.Cells.Item(107, 9).Value = 10206195
$ws.Cells.Item(107, 10).Value = 4000
$ws.Cells.Item(107, 11).Value = 10206195
$ws.Cells.Item(107, 12).Value = 4000
$ws.Cells.Item(107, 13).Value = -10204275
$ws.Cells.Item(107, 14).Value = -7840

$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 13).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 286
$ws.Cells.Item(22, 9).Value = 286
$ws.Cells.Item(22, 11).Value = 286
$ws.Cells.Item(22, 13).Value = 64

$ws.Cells.Item(58, 8).Value = 2393.8076
$ws.Cells.Item(58, 10).Value = 3404.1
$ws.Cells.Item(58, 12).Value = 3404.1
$ws.Cells.Item(58, 14).Value = -3810.1

$ws.Cells.Item(99, 8).Value = 3658.1667
$ws.Cells.Item(99, 9).Value = 3251.375
$ws.Cells.Item(99, 11).Value = 3251.375
$ws.Cells.Item(99, 13).Value = -1753.375

$ws.Cells.Item(103, 8).Value = 130000
$ws.Cells.Item(103, 9).Value = 0
$ws.Cells.Item(103, 11).Value = 0
$ws.Cells.Item(103, 13).ClearContents()

$ws.Cells.Item(105, 8).Value = 2597.5
$ws.Cells.Item(105, 9).Value = 2195
$ws.Cells.Item(105, 11).Value = 2195
$ws.Cells.Item(105, 13).Value = -448

$ws.Cells.Item(126, 8).Value = 3658.1667
$ws.Cells.Item(126, 9).Value = 3251.375
$ws.Cells.Item(126, 11).Value = 9754.125
$ws.Cells.Item(126, 13).Value = -7284.125

$ws.Cells.Item(134, 8).Value = 3561.2307
$ws.Cells.Item(134, 9).Value = 3183.625
$ws.Cells.Item(134, 11).Value = 9550.875
$ws.Cells.Item(134, 13).Value = -7015.875

$ws.Cells.Item(136, 8).Value = 2393.8076
$ws.Cells.Item(136, 10).Value = 3404.1
$ws.Cells.Item(136, 12).Value = 10212.3
$ws.Cells.Item(136, 14).Value = -15312.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 19843610
$ws.Cells.Item(131, 10).Value = 25644374
$ws.Cells.Item(131, 12).Value = 76933122
$ws.Cells.Item(131, 14).Value = -76943202

$ws.Cells.Item(139, 8).Value = 62502228
$ws.Cells.Item(139, 9).Value = 83334810
$ws.Cells.Item(139, 11).Value = 250004430
$ws.Cells.Item(139, 13).Value = -249999290

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 6905.0625
$ws.Cells.Item(2, 9).Value = 939.9091
$ws.Cells.Item(2, 10).Value = 20028.4
$ws.Cells.Item(2, 11).Value = 939.9091
$ws.Cells.Item(2, 12).Value = 20028.4
$ws.Cells.Item(2, 13).Value = -826.9091
$ws.Cells.Item(2, 14).Value = -20254.4

$ws.Cells.Item(70, 8).Value = 20008238
$ws.Cells.Item(70, 9).Value = 25008922
$ws.Cells.Item(70, 11).Value = 25008922
$ws.Cells.Item(70, 13).Value = -25008652

$ws.Cells.Item(73, 8).Value = 20008238
$ws.Cells.Item(73, 9).Value = 25008922
$ws.Cells.Item(73, 11).Value = 25008922
$ws.Cells.Item(73, 13).Value = -25007986

$ws.Cells.Item(102, 8).Value = 7498376.5
$ws.Cells.Item(102, 9).Value = 7938590
$ws.Cells.Item(102, 11).Value = 7938590
$ws.Cells.Item(102, 13).Value = -7936968

$ws.Cells.Item(122, 8).Value = 448841.9
$ws.Cells.Item(122, 9).Value = 559301.9399999999
$ws.Cells.Item(122, 11).Value = 1677905.82
$ws.Cells.Item(122, 13).Value = -1675455.82

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 37040308
$ws.Cells.Item(93, 9).Value = 41669972
$ws.Cells.Item(93, 11).Value = 41669972
$ws.Cells.Item(93, 13).Value = -41668724

$ws.Cells.Item(123, 8).Value = 39999
$ws.Cells.Item(123, 10).Value = 39999
$ws.Cells.Item(123, 12).Value = 39999
$ws.Cells.Item(123, 14).Value = -49799

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(70, 8).Value = 60000
$ws.Cells.Item(70, 10).Value = 60000
$ws.Cells.Item(70, 12).Value = 60000
$ws.Cells.Item(70, 14).Value = -60630

$ws.Cells.Item(73, 8).Value = 60000
$ws.Cells.Item(73, 10).Value = 60000
$ws.Cells.Item(73, 12).Value = 60000
$ws.Cells.Item(73, 14).Value = -62184

$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 13).ClearContents()

$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 13).ClearContents()

$ws.Cells.Item(126, 8).Value = 2579.0356
$ws.Cells.Item(126, 9).Value = 2518.1924
$ws.Cells.Item(126, 11).Value = 7554.5772
$ws.Cells.Item(126, 13).Value = -5084.5772

$ws.Cells.Item(132, 8).Value = 20021312
$ws.Cells.Item(132, 9).Value = 34487748
$ws.Cells.Item(132, 10).Value = 43851.953
$ws.Cells.Item(132, 11).Value = 103463244
$ws.Cells.Item(132, 12).Value = 131555.859
$ws.Cells.Item(132, 13).Value = -103460714
$ws.Cells.Item(132, 14).Value = -136615.859

$ws.Cells.Item(136, 8).Value = 2180.6572
$ws.Cells.Item(136, 9).Value = 1353.16
$ws.Cells.Item(136, 10).Value = 4249.4
$ws.Cells.Item(136, 11).Value = 4059.48
$ws.Cells.Item(136, 12).Value = 12748.2
$ws.Cells.Item(136, 13).Value = -1509.48
$ws.Cells.Item(136, 14).Value = -17848.2
